# Screen and Transition Logic Implementation
# Move the "Weapons" (J) and "Items" (K) columns from the "Resources" sheet
# to the "Tools" sheet (as columns B and C), and tweak the active-cell
# selection on the "Tech Tree" sheet.

$wb = $excel.ActiveWorkbook

$wsResources = $wb.Worksheets.Item("Resources")
$wsTools     = $wb.Worksheets.Item("Tools")
$wsTechTree  = $wb.Worksheets.Item("Tech Tree")

# --- Move data from Resources!J2:K11 -> Tools!B2:C11 ---
$srcRange = $wsResources.Range("J2:K11")
$dstRange = $wsTools.Range("B2:C11")
$dstRange.Value2 = $srcRange.Value2

# Copy the style (bold+underline header) of the header row J2:K2 -> B2:C2
$wsResources.Range("J2:K2").Copy()
$wsTools.Range("B2:C2").PasteSpecial(-4122)  # xlPasteFormats

# Clear the now-empty source columns on Resources
$srcRange.Clear()

# --- Column widths ---
# Resources: columns J & K shrink to the sheet's default (no longer bestFit)
$wsResources.Columns.Item(10).ColumnWidth = 8.32
$wsResources.Columns.Item(11).ColumnWidth = 8.32

# Tools: columns B & C take on the bestFit-like widths that used to belong to J & K
$wsTools.Columns.Item(2).ColumnWidth = 14.1
$wsTools.Columns.Item(3).ColumnWidth = 9

# --- Selections / active cells ---
$wsTools.Range("D34").Select()
$wsTechTree.Range("D2").Select()
$wsResources.Range("H18").Select()

$wb.Save()
